# Edit CovidUK.xlsx: append newest day (2021-01-27), rename sheet, and
# relocate backfilled rows (dates 2020-10-20 / 2020-09-05) that the source
# feed re-supplied out of chronological order to the bottom of the sheet,
# while the previously-misplaced 2020-06-14 row settles back into its
# correct sorted position. Net effect vs. the original file: one new row
# (392) is added and 184 existing rows shift to absorb it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet to reflect the new "as of" date.
$ws.Name = "overview_2021-01-27"

# 2. Full target state for A/E/F across rows 2..392 (date serial, newCases,
#    newDeaths), as "date,cases,deaths" records separated by ";".
$data = "44223,25308,1725;44222,20089,1631;44221,22195,592;44220,30004,610;44219,33552,1348;44218,40261,1401;44217,37892,1290;44216,38905,1820;44215,33355,1610;44214,37535,599;44213,38598,671;44212,41346,1295;44211,55761,1280;44210,48682,1248;44209,47525,1564;44208,45533,1243;44207,46169,529;44206,54940,563;44205,59937,1035;44204,68053,1325;44203,52618,1162;44202,62322,1041;44201,60916,830;44200,58784,407;44199,54990,455;44198,57725,445;44197,53285,656;44196,55892,964;44195,50023,981;44194,53135,414;44193,41385,357;44192,30501,317;44191,35691,230;44190,32725,613;44189,39877,585;44188,39237,744;44187,36804,691;44186,33364,215;44185,35928,326;44184,27052,534;44183,28507,489;44182,35383,532;44181,25161,613;44180,18450,506;44179,20263,232;44178,18447,144;44177,21502,519;44176,21672,424;44175,20964,516;44174,16578,533;44173,12282,616;44172,14718,189;44171,17272,231;44170,15539,397;44169,16298,504;44168,14879,414;44167,16170,648;44166,13430,603;44165,12330,205;44164,12155,215;44163,15871,479;44162,16022,521;44161,17555,498;44160,18213,696;44159,11299,608;44158,15450,206;44157,18662,398;44156,19875,341;44155,20252,511;44154,22915,501;44153,19609,529;44152,20051,598;44151,21363,213;44150,24962,168;44149,26860,462;44148,27301,376;44147,33470,563;44146,22950,595;44145,20412,532;44144,21350,194;44143,20572,156;44142,24957,413;44141,23287,355;44140,24141,378;44139,25177,492;44138,20018,397;44137,18950,136;44136,23254,162;44135,21915,326;44134,24405,274;44133,23065,280;44132,24701,310;44131,22885,367;44130,20890,102;44129,19790,151;44128,23012,174;44127,20530,224;44126,21242,189;44125,26688,191;44123,18804,80;44122,16982,67;44121,16171,150;44120,15650,136;44119,18980,138;44118,19724,137;44117,17234,143;44116,13972,50;44115,12872,65;44114,15166,81;44113,13864,87;44112,17540,77;44111,14162,70;44110,14542,76;44109,12594,19;44108,22961,33;44107,12872,49;44106,6968,66;44105,6914,59;44104,7108,71;44103,7143,71;44102,4044,13;44101,5693,17;44100,6042,34;44099,6874,35;44098,6634,40;44097,6178,37;44096,4926,37;44095,4368,11;44094,3899,18;44093,4422,27;44092,4322,27;44091,3395,21;44090,3991,20;44089,3105,27;44088,2621,9;44087,3330,5;44086,3497,9;44085,3539,6;44084,2919,14;44083,2659,8;44082,2460,32;44081,2948,3;44080,2988,2;44078,1940,10;44077,1735,13;44076,1508,10;44075,1295,3;44074,1406,2;44073,1715,1;44072,1108,12;44071,1276,9;44070,1522,12;44069,1048,16;44068,1184,16;44067,853,4;44066,1160,6;44065,1288,18;44064,1033,2;44063,1182,6;44062,812,16;44061,1089,12;44060,713,3;44059,1040,5;44058,1077,3;44057,1441,11;44056,1129,18;44055,1009,20;44054,1148,13;44053,740,18;44052,1062,5;44051,758,3;44050,871,12;44049,950,18;44048,892,14;44047,670,18;44046,898,1;44045,744,5;44044,761,13;44043,880,20;44042,846,0;44041,763,34;44040,581,21;44039,685,3;44038,745,8;44037,767,15;44036,768,32;44035,769,9;44034,560,17;44033,445,25;44032,557,10;44031,726,11;44030,827,9;44029,687,26;44028,642,24;44027,538,26;44026,398,44;44025,530,10;44024,650,9;44023,820,17;44022,512,34;44021,642,31;44020,630,57;44019,581,54;44018,343,11;44017,516,19;44016,624,32;44015,544,50;44014,576,41;44013,829,97;44012,689,53;44011,815,21;44010,901,31;44009,890,40;44008,1006,79;44007,1118,99;44006,653,87;44005,874,94;44004,958,14;44003,1221,31;44002,1295,71;44001,1346,84;44000,1218,67;43999,1115,110;43998,1279,120;43997,1056,29;43996,1514,27;43995,1425,107;43994,1541,131;43993,1266,76;43992,1003,164;43991,1387,195;43990,1205,47;43989,1326,54;43988,1557,143;43987,1650,258;43986,1805,130;43985,1871,254;43984,1613,249;43983,1570,86;43982,1936,60;43981,2445,154;43980,2095,274;43979,1887,343;43978,2013,422;43977,2004,131;43976,1625,104;43975,2409,379;43974,2959,220;43973,3287,291;43972,2615,273;43971,2472,328;43970,2412,500;43969,2684,146;43968,3562,67;43967,3451,411;43966,3560,350;43965,3446,352;43964,3242,447;43963,3403,614;43962,3877,187;43961,3923,217;43960,3896,275;43959,4649,579;43958,5614,458;43957,6111,647;43956,4406,726;43955,3985,272;43954,4339,253;43953,4806,584;43952,6201,698;43951,6032,634;43950,4076,769;43949,3996,969;43948,4310,320;43947,4463,364;43946,4913,815;43945,5386,1010;43944,4583,682;43943,4451,847;43942,4301,1224;43941,4676,570;43940,5850,432;43939,5526,1105;43938,5599,913;43937,4618,1036;43936,4605,880;43935,5252,1076;43934,4342,724;43933,5288,657;43932,5234,843;43931,5706,1122;43930,4675,1116;43929,5865,1030;43928,3888,1105;43927,4143,567;43926,6199,599;43925,4000,756;43924,4672,736;43923,4522,657;43922,4567,672;43921,3250,403;43920,2665,374;43919,2502,212;43918,2556,292;43917,2890,288;43916,2129,181;43915,1452,191;43914,1427,148;43913,967,76;43912,665,36;43911,1035,58;43910,714,32;43909,643,45;43908,676,34;43907,407,16;43906,152,22;43905,330,14;43904,264,18;43903,207,1;43902,134,2;43901,83,0;43900,52,4;43899,50,1;43898,65,0;43897,46,1;43896,46,1;43895,29,0;43894,34,0;43893,11,0;43892,5,0;43891,12,0;43890,4,0;43889,6,0;43888,0,0;43887,0,0;43886,0,0;43885,4,0;43884,0,0;43883,0,0;43882,0,0;43881,0,0;43880,0,0;43879,0,0;43878,0,0;43877,0,0;43876,0,0;43875,0,0;43874,1,0;43873,0,0;43872,0,0;43871,4,0;43870,1,0;43869,0,0;43868,0,0;43867,1,0;43866,0,0;43865,0,0;43864,0,0;43863,0,0;43862,0,0;43861,2,0;43860,0,0;43859,0,0;43858,0,0;43857,0,0;43856,0,0;43855,0,0;43854,0,0;43853,0,0;43852,0,0;43851,0,0;43850,0,0;43849,0,0;43848,0,0;43847,0,0;43846,0,0;43845,0,0;43844,0,0;43843,0,0;43842,0,0;43841,0,0;43840,0,0;43839,0,0;43838,0,0;43837,0,0;43836,0,0;43835,0,0;43834,0,0;43833,0,0;44124,21331,241;44079,1813,12"
$records = $data -split ";"

for ($i = 0; $i -lt $records.Count; $i++) {
    $parts = $records[$i] -split ","
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = [double]$parts[0]
    $ws.Cells.Item($row, 5).Value = [double]$parts[1]
    $ws.Cells.Item($row, 6).Value = [double]$parts[2]
}

# 3. Row 392 is brand new -- give it the same areaType/areaCode/areaName
#    text as every other data row, and copy the date number format from an
#    existing, untouched date cell onto its A cell.
$ws.Cells.Item(392, 2).Value = "overview"
$ws.Cells.Item(392, 3).Value = "K02000001"
$ws.Cells.Item(392, 4).Value = "United Kingdom"

$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(392, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

"done"
